$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.832.98"
$ws.Cells.Item(2, 5).Value = "  -1.55%  "

$ws.Cells.Item(3, 4).Value = "3.496.53"
$ws.Cells.Item(3, 5).Value = "  -0.69%  "

$ws.Cells.Item(4, 5).Value = "  -0.17%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "601.72"
$ws.Cells.Item(5, 5).Value = "  -1.41%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "148.05"
$ws.Cells.Item(6, 5).Value = "  -3.00%  "

$ws.Cells.Item(7, 4).Value = "3.495.20"
$ws.Cells.Item(7, 5).Value = "  -0.61%  "

$ws.Cells.Item(8, 5).Value = "  -0.01%  "

$ws.Cells.Item(9, 5).Value = "  -1.85%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.143"
$ws.Cells.Item(10, 5).Value = "  -0.99%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "7.99"
$ws.Cells.Item(11, 5).Value = "  +5.23%  "

$ws.Cells.Item(12, 5).Value = "  -2.63%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000214"
$ws.Cells.Item(13, 5).Value = "  -1.71%  "

$ws.Cells.Item(14, 4).Value = "4.087.93"
$ws.Cells.Item(14, 5).Value = "  -0.77%  "

$ws.Cells.Item(15, 5).Value = "  -5.11%  "

$ws.Cells.Item(16, 4).Value = "3.492.39"
$ws.Cells.Item(16, 5).Value = "  -1.00%  "

$ws.Cells.Item(17, 4).Value = "66.817.99"
$ws.Cells.Item(17, 5).Value = "  -1.46%  "

$ws.Cells.Item(18, 5).Value = "  +0.27%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.39"
$ws.Cells.Item(19, 5).Value = "  -3.40%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "10.37"
$ws.Cells.Item(20, 5).Value = "  +5.62%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "15.32"
$ws.Cells.Item(21, 5).Value = "  -2.22%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "434.08"
$ws.Cells.Item(22, 5).Value = "  -3.75%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.610"
$ws.Cells.Item(23, 5).Value = "  -4.11%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "79.59"
$ws.Cells.Item(24, 5).Value = "  +1.93%  "

$ws.Cells.Item(25, 4).Value = "3.633.21"
$ws.Cells.Item(25, 5).Value = "  -0.99%  "

$ws.Cells.Item(26, 5).Value = "  +0.04%  "

$ws.Cells.Item(27, 5).Value = "  -7.55%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "9.83"
$ws.Cells.Item(28, 5).Value = "  -3.21%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "8.24"
$ws.Cells.Item(29, 5).Value = "  -8.71%  "

$ws.Cells.Item(30, 5).Value = "  -1.18%  "

$ws.Cells.Item(31, 5).Value = "  -4.26%  "

$ws.Cells.Item(32, 5).Value = "  +0.01%  "

$ws.Cells.Item(33, 5).Value = "  -3.57%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "25.41"
$ws.Cells.Item(34, 5).Value = "  -1.70%  "

$ws.Cells.Item(35, 4).Value = "3.489.26"
$ws.Cells.Item(35, 5).Value = "  -0.74%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "5.96"
$ws.Cells.Item(36, 5).Value = "  -4.49%  "

$ws.Cells.Item(37, 5).Value = "  -4.49%  "

$ws.Cells.Item(38, 5).Value = "  -1.17%  "

$ws.Cells.Item(39, 5).Value = "  +0.00%  "

$ws.Cells.Item(40, 5).Value = "  -0.21%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0892"
$ws.Cells.Item(41, 5).Value = "  -1.48%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "170.62"
$ws.Cells.Item(42, 5).Value = "  -2.06%  "

$ws.Cells.Item(43, 2).Value = "Stacks"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.08"
$ws.Cells.Item(43, 5).Value = "  -10.81%  "

$ws.Cells.Item(44, 2).Value = "Filecoin"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.42"
$ws.Cells.Item(44, 5).Value = "  -2.87%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.896"
$ws.Cells.Item(45, 5).Value = "  +1.29%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "28.95"
$ws.Cells.Item(46, 5).Value = "  -6.15%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "45.81"
$ws.Cells.Item(47, 5).Value = "  -2.01%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.24"
$ws.Cells.Item(48, 5).Value = "  -6.26%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.48"
$ws.Cells.Item(49, 5).Value = "  -2.96%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.42"
$ws.Cells.Item(50, 5).Value = "  -5.49%  "

$ws.Cells.Item(51, 2).Value = "TheGraph"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.247"
$ws.Cells.Item(51, 5).Value = "  -3.57%  "
